$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "TODO Before 0.0.1" sheet: mark a handful of crystal/gem related
#    backlog items as "in-progress" (were "todo").
# ------------------------------------------------------------------
$todoSheet = $wb.Worksheets.Item("TODO Before 0.0.1")
$todoSheet.Range("C33").Value = "in-progress"
$todoSheet.Range("C34").Value = "in-progress"
$todoSheet.Range("C35").Value = "in-progress"
$todoSheet.Range("C37").Value = "in-progress"
$todoSheet.Range("C44").Value = "in-progress"

# ------------------------------------------------------------------
# 2. "Logs" sheet: append a new dev-log entry for today's work.
# ------------------------------------------------------------------
$logsSheet = $wb.Worksheets.Item("Logs")

# Copy the formatting of the last existing row down onto the new row
# so the new cells inherit the same date / wrap-text styling.
$logsSheet.Range("A64:B64").Copy() | Out-Null
$logsSheet.Range("A65:B65").PasteSpecial(-4122) | Out-Null

$logsSheet.Range("A65").Value = 45613
$logsSheet.Range("B65").Value = "add crystal, crystal controller item/skill, fix fontain"

# ------------------------------------------------------------------
# 3. Restore on-screen selections to match where the author left off.
# ------------------------------------------------------------------
$todoSheet.Activate() | Out-Null
$todoSheet.Range("C33").Select() | Out-Null

$logsSheet.Activate() | Out-Null
$logsSheet.Range("B66").Select() | Out-Null
